$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: bump the "Latest HO Xliff Generate Date" for the
#     542fba49... report (row 2) and the 5bdf83ce... report (row 3) to the
#     freshly (re)generated timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 09:55:59"
$wsOverview.Range("G3").Value = "2016-08-22 09:55:59"

# --- "zh-cn" sheet: handback type changed from "ht" to "mt" and the
#     handoff/handback datetimes for the 542fba49... file were refreshed.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-22 09:55:55"
$wsZhCn.Range("H3").Value = "2016-08-22 09:55:55"
$wsZhCn.Range("K2").Value = "2016-08-22 09:56:19"
$wsZhCn.Range("K3").Value = "2016-08-22 09:56:19"

# --- "de-de" sheet: same handback-type change, plus the handoff datetime
#     (shared with Overview's G column string) and handback datetime refresh.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-22 09:55:59"
$wsDeDe.Range("H3").Value = "2016-08-22 09:55:59"
$wsDeDe.Range("K2").Value = "2016-08-22 09:56:26"
$wsDeDe.Range("K3").Value = "2016-08-22 09:56:26"
